$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.335.97'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '3.800.86'
$ws.Range("E3").Value = '  +3.40%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '619.06'
$ws.Range("E5").Value = '  +3.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.17'
$ws.Range("E6").Value = '  -4.44%  '
$ws.Range("D7").Value = '3.800.92'
$ws.Range("E7").Value = '  +3.48%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.537'
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  +3.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.34'
$ws.Range("E11").Value = '  -3.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.495'
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.07'
$ws.Range("E13").Value = '  +2.80%  '
$ws.Range("E14").Value = '  +0.05%  '
$ws.Range("D15").Value = '4.430.08'
$ws.Range("E15").Value = '  +3.26%  '
$ws.Range("D16").Value = '3.804.30'
$ws.Range("E16").Value = '  +4.03%  '
$ws.Range("D17").Value = '70.352.50'
$ws.Range("E17").Value = '  -0.83%  '
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '515.87'
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.69'
$ws.Range("E21").Value = '  -3.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.64'
$ws.Range("E22").Value = '  +3.55%  '
$ws.Range("E23").Value = '  -3.43%  '
$ws.Range("E24").Value = '  +3.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.22'
$ws.Range("E25").Value = '  -0.59%  '
$ws.Range("E26").Value = '  -2.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.14'
$ws.Range("E27").Value = '  +2.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000138'
$ws.Range("E28").Value = '  +24.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("E32").Value = '  -5.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.86'
$ws.Range("E33").Value = '  -1.06%  '
$ws.Range("E34").Value = '  -2.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("E37").Value = '  +2.25%  '
$ws.Range("E38").Value = '  +0.68%  '
$ws.Range("E39").Value = '  +1.86%  '
$ws.Range("E40").Value = '  +2.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.12'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("B42").Value = 'Cosmos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.78'
$ws.Range("E42").Value = '  -1.85%  '
$ws.Range("B43").Value = 'Arweave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '44.12'
$ws.Range("E43").Value = '  -7.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '423.05'
$ws.Range("E44").Value = '  +3.76%  '
$ws.Range("D45").Value = '3.068.97'
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.76'
$ws.Range("E46").Value = '  -2.06%  '
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.60'
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '135.83'
$ws.Range("E49").Value = '  +0.74%  '
$ws.Range("E51").Value = '  -0.36%  '
